$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.286.40"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "'  -2.91%  "
$ws.Range("E2").Style = $ws.Range("B2").Style
$ws.Range("D3").Value = "'1.554.17"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "'  -4.69%  "
$ws.Range("E3").Style = $ws.Range("B3").Style
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = $ws.Range("B4").Style
$ws.Range("D5").Value = "'207.15"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "'  -3.35%  "
$ws.Range("E5").Style = $ws.Range("B5").Style
$ws.Range("E6").Value = "'  -0.04%  "
$ws.Range("E6").Style = $ws.Range("B6").Style
$ws.Range("D7").Value = "'0.477"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "'  -5.26%  "
$ws.Range("E7").Style = $ws.Range("B7").Style
$ws.Range("E8").Value = "'  -1.58%  "
$ws.Range("E8").Style = $ws.Range("B8").Style
$ws.Range("E9").Value = "'  -3.21%  "
$ws.Range("E9").Style = $ws.Range("B9").Style
$ws.Range("D10").Value = "'17.77"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "'  -4.62%  "
$ws.Range("E10").Style = $ws.Range("B10").Style
$ws.Range("E11").Value = "'  -0.99%  "
$ws.Range("E11").Style = $ws.Range("B11").Style
$ws.Range("D12").Value = "'1.768.08"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "'  -4.80%  "
$ws.Range("E12").Style = $ws.Range("B12").Style
$ws.Range("D13").Value = "'1.547.60"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "'  -4.74%  "
$ws.Range("E13").Style = $ws.Range("B13").Style
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "'  -4.58%  "
$ws.Range("E14").Style = $ws.Range("B14").Style
$ws.Range("E15").Value = "'  -4.44%  "
$ws.Range("E15").Style = $ws.Range("B15").Style
$ws.Range("D16").Value = "'25.271.49"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E17").Value = "'  -4.77%  "
$ws.Range("E17").Style = $ws.Range("B17").Style
$ws.Range("D18").Value = "'58.77"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "'  -4.67%  "
$ws.Range("E18").Style = $ws.Range("B18").Style
$ws.Range("E19").Value = "'  +0.01%  "
$ws.Range("E19").Style = $ws.Range("B19").Style
$ws.Range("D20").Value = "'185.69"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "'  -3.82%  "
$ws.Range("E20").Style = $ws.Range("B20").Style
$ws.Range("E21").Value = "'  -3.60%  "
$ws.Range("E21").Style = $ws.Range("B21").Style
$ws.Range("D22").Value = "'9.28"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "'  -3.02%  "
$ws.Range("E22").Style = $ws.Range("B22").Style
$ws.Range("E23").Value = "'  -3.71%  "
$ws.Range("E23").Style = $ws.Range("B23").Style
$ws.Range("E24").Value = "'  -4.20%  "
$ws.Range("E24").Style = $ws.Range("B24").Style
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "'  -0.06%  "
$ws.Range("E25").Style = $ws.Range("B25").Style
$ws.Range("D26").Value = "'139.86"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "'  -3.07%  "
$ws.Range("E26").Style = $ws.Range("B26").Style
$ws.Range("E27").Value = "'  -5.06%  "
$ws.Range("E27").Style = $ws.Range("B27").Style
$ws.Range("D28").Value = "'14.88"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "'  -2.63%  "
$ws.Range("E28").Style = $ws.Range("B28").Style
$ws.Range("D29").Value = "'6.40"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "'  -5.03%  "
$ws.Range("E29").Style = $ws.Range("B29").Style
$ws.Range("E30").Value = "'  -6.92%  "
$ws.Range("E30").Style = $ws.Range("B30").Style
$ws.Range("D31").Value = "'0.0467"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "'  -3.44%  "
$ws.Range("E31").Style = $ws.Range("B31").Style
$ws.Range("E32").Value = "'  -3.42%  "
$ws.Range("E32").Style = $ws.Range("B32").Style
$ws.Range("E33").Value = "'  -5.03%  "
$ws.Range("E33").Style = $ws.Range("B33").Style
$ws.Range("E34").Value = "'  -3.34%  "
$ws.Range("E34").Style = $ws.Range("B34").Style
$ws.Range("E35").Value = "'  -3.44%  "
$ws.Range("E35").Style = $ws.Range("B35").Style
$ws.Range("D36").Value = "'1.085.16"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "'  -3.98%  "
$ws.Range("E36").Style = $ws.Range("B36").Style
$ws.Range("E37").Value = "'  -0.08%  "
$ws.Range("E37").Style = $ws.Range("B37").Style
$ws.Range("E38").Value = "'  -2.92%  "
$ws.Range("E38").Style = $ws.Range("B38").Style
$ws.Range("E39").Value = "'  -4.90%  "
$ws.Range("E39").Style = $ws.Range("B39").Style
$ws.Range("D40").Value = "'0.767"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "'  -10.22%  "
$ws.Range("E40").Style = $ws.Range("B40").Style
$ws.Range("E41").Value = "'  -7.83%  "
$ws.Range("E41").Style = $ws.Range("B41").Style
$ws.Range("D42").Value = "'0.799"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "'  +5.14%  "
$ws.Range("E42").Style = $ws.Range("B42").Style
$ws.Range("D43").Value = "'92.59"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "'  -5.86%  "
$ws.Range("E43").Style = $ws.Range("B43").Style
$ws.Range("D44").Value = "'5.05"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "'  -1.69%  "
$ws.Range("E44").Style = $ws.Range("B44").Style
$ws.Range("D45").Value = "'1.682.75"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "'  -4.77%  "
$ws.Range("E45").Style = $ws.Range("B45").Style
$ws.Range("E46").Value = "'  -2.62%  "
$ws.Range("E46").Style = $ws.Range("B46").Style
$ws.Range("E47").Value = "'  -1.94%  "
$ws.Range("E47").Style = $ws.Range("B47").Style
$ws.Range("D48").Value = "'52.36"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "'  -4.00%  "
$ws.Range("E48").Style = $ws.Range("B48").Style
$ws.Range("E49").Value = "'  -4.03%  "
$ws.Range("E49").Style = $ws.Range("B49").Style
$ws.Range("E50").Value = "'  -0.25%  "
$ws.Range("E50").Style = $ws.Range("B50").Style
$ws.Range("E51").Value = "'  -2.19%  "
$ws.Range("E51").Style = $ws.Range("B51").Style
